$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A1").Value = 2619.6710337398545
$ws.Range("B1").Value = 1744.3261301353346
$ws.Range("C1").Value = 1744.3746610075209

$ws.Range("A2").Value = 2422.980424162788
$ws.Range("B2").Value = 1403.7753353357568
$ws.Range("C2").Value = 1347.2865867830023

$ws.Range("A3").Value = 2827.7480600865219
$ws.Range("B3").Value = 2266.966322519324
$ws.Range("C3").Value = 1933.4995920571732

$ws.Range("A4").Value = 2539.7595924399211
$ws.Range("B4").Value = 1813.7919403213787
$ws.Range("C4").Value = 1648.8017994991731

$ws.Range("A5").Value = 2844.8584118708213
$ws.Range("B5").Value = 1976.0549125908847
$ws.Range("C5").Value = 2168.3385675552581

$ws.Range("A6").Value = 2552.4329862565151
$ws.Range("B6").Value = 1795.1157047345025
$ws.Range("C6").Value = 1891.7084305962494

$ws.Range("A7").Value = 2396.6531084939179
$ws.Range("B7").Value = 1907.6776304748842
$ws.Range("C7").Value = 1667.4837984331457

$ws.Range("A8").Value = 2988.4021092342487
$ws.Range("B8").Value = 2463.1835649003879
$ws.Range("C8").Value = 2321.1985563551616

$ws.Range("A10").Value = 2442.5594633529463
$ws.Range("B10").Value = 1491.3792481496946
$ws.Range("C10").Value = 1415.3482452983621

$ws.Range("A11").Value = 2285.3777130740923
$ws.Range("B11").Value = 1678.0274308615462
$ws.Range("C11").Value = 1620.9433178454256

$ws.Range("A12").Value = 3107.6240713388847
$ws.Range("B12").Value = 2479.7036053696779
$ws.Range("C12").Value = 2361.6320543880497
